$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("1000 Bs = 6.13 = 25026.84 pesos", "1000 Bs = 6.12 = 24770.49 pesos")
$text = $text.Replace("25026.84 pesos = 6.14 = 982.22 Bs", "24770.49 pesos = 6.08 = 962.24 Bs")
$cell.Value = $text

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 163.501
$wsTasas.Range("O10").Value = 4050
$wsTasas.Range("N12").Value = 4077
$wsTasas.Range("O12").Value = 158.376
